$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 10

$ws.Cells.Item($row, 1).Value = "2025-08-01 06:26:44"
$ws.Cells.Item($row, 2).Value = "create-team"
$ws.Cells.Item($row, 3).Value = "new-organization97"
$ws.Cells.Item($row, 4).Value = "firstteam"
$ws.Cells.Item($row, 5).Value = "demo"
$ws.Cells.Item($row, 6).Value = "Vignesh2122"
$ws.Cells.Item($row, 7).Value = "pull"

# "False" would otherwise be auto-coerced to a Boolean by Excel's type
# inference; force it to remain plain text like the other string cells,
# then clear the resulting quote-prefix style so no extra style index
# is introduced on the cell.
$ws.Cells.Item($row, 9).Value = "'False"
$ws.Cells.Item($row, 9).Style = "Normal"
